$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.925.81'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '2.502.39'
$ws.Range('E3').Value = '  +2.39%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '539.72'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.62'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.12%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.570'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('D9').Value = '2.526.22'
$ws.Range('E9').Value = '  +2.74%  '
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('E11').Value = '  +0.68%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.61'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +6.90%  '
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').Value = '2.943.35'
$ws.Range('E14').Value = '  +2.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.49'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('D16').Value = '58.861.89'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000140'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.01%  '
$ws.Range('D18').Value = '2.507.38'
$ws.Range('E18').Value = '  +0.52%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.20'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.27'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.59'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.74%  '
$ws.Range('E22').Value = '  +3.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.78'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.96'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.438'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.66%  '
$ws.Range('E26').Value = '  +1.08%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.614.81'
$ws.Range('E27').Value = '  +2.67%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.995'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.76'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.22%  '
$ws.Range('D30').Value = '0.0₃0773'
$ws.Range('E30').Value = '  +0.72%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.67'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.52%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.80'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('E33').Value = '  -4.11%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '156.76'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.43'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +3.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.69'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.34'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.60'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -7.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.66'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.84'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '297.14'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.70'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.824'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.601'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.44%  '
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0929'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.51'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.65'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.13%  '
$ws.Range('E51').Value = '  +0.24%  '
